# Scheduled-runner update: refresh market-price derived columns
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
#  LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ, columns H:N)
# on the Kujata_Profits workbook's per-job sheets (ALC, ARM, BSM, CRP,
# CUL, LTW, WVR) for a handful of leves whose underlying item prices moved.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H17").Value = 1580
$ws_ALC.Range("J17").Value = 1821.4286
$ws_ALC.Range("L17").Value = 5464.2858
$ws_ALC.Range("N17").Value = -5800.2858
$ws_ALC.Range("H74").Value = 3167
$ws_ALC.Range("I74").Value = 4003
$ws_ALC.Range("J74").Value = 2999.8
$ws_ALC.Range("K74").Value = 4003
$ws_ALC.Range("L74").Value = 2999.8
$ws_ALC.Range("M74").Value = -3067
$ws_ALC.Range("N74").Value = -4871.8
$ws_ALC.Range("H77").Value = 3167
$ws_ALC.Range("I77").Value = 4003
$ws_ALC.Range("J77").Value = 2999.8
$ws_ALC.Range("K77").Value = 20015
$ws_ALC.Range("L77").Value = 14999
$ws_ALC.Range("M77").Value = -15335
$ws_ALC.Range("N77").Value = -24359
$ws_ALC.Range("H116").Value = 3599.8462
$ws_ALC.Range("I116").Value = 2970.5715
$ws_ALC.Range("K116").Value = 2970.5715
$ws_ALC.Range("M116").Value = 471.4285
$ws_ALC.Range("H133").Value = 34459
$ws_ALC.Range("J133").Value = 34459
$ws_ALC.Range("L133").Value = 34459
$ws_ALC.Range("N133").Value = -44579

# ---- ARM ----
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H2").Value = 907.8077
$ws_ARM.Range("I2").Value = 724.7895
$ws_ARM.Range("K2").Value = 724.7895
$ws_ARM.Range("M2").Value = -611.7895
$ws_ARM.Range("H13").Value = 9600400
$ws_ARM.Range("I13").Value = 12000250
$ws_ARM.Range("J13").Value = 999
$ws_ARM.Range("K13").Value = 12000250
$ws_ARM.Range("L13").Value = 999
$ws_ARM.Range("M13").Value = -12000106
$ws_ARM.Range("N13").Value = -1287
$ws_ARM.Range("H32").Value = 7249.1226
$ws_ARM.Range("I32").Value = 5956.381
$ws_ARM.Range("J32").Value = 15005.571
$ws_ARM.Range("K32").Value = 5956.381
$ws_ARM.Range("L32").Value = 15005.571
$ws_ARM.Range("M32").Value = -5669.381
$ws_ARM.Range("N32").Value = -15579.571
$ws_ARM.Range("H63").Value = 29414420
$ws_ARM.Range("I63").Value = 2000.1818
$ws_ARM.Range("J63").Value = 83337190
$ws_ARM.Range("K63").Value = 2000.1818
$ws_ARM.Range("L63").Value = 83337190
$ws_ARM.Range("M63").Value = -1314.1818
$ws_ARM.Range("N63").Value = -83338562
$ws_ARM.Range("H66").Value = 29414420
$ws_ARM.Range("I66").Value = 2000.1818
$ws_ARM.Range("J66").Value = 83337190
$ws_ARM.Range("K66").Value = 10000.909
$ws_ARM.Range("L66").Value = 416685950
$ws_ARM.Range("M66").Value = -6568.909
$ws_ARM.Range("N66").Value = -416692814
$ws_ARM.Range("H116").Value = 907.8077
$ws_ARM.Range("I116").Value = 724.7895
$ws_ARM.Range("K116").Value = 724.7895
$ws_ARM.Range("M116").Value = 1569.2105
$ws_ARM.Range("H132").Value = 2825.8647
$ws_ARM.Range("I132").Value = 2263.318
$ws_ARM.Range("J132").Value = 3650.9333
$ws_ARM.Range("K132").Value = 6789.954000000001
$ws_ARM.Range("L132").Value = 10952.7999
$ws_ARM.Range("M132").Value = -4259.954000000001
$ws_ARM.Range("N132").Value = -16012.7999

# ---- BSM ----
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H3").Value = 907.8077
$ws_BSM.Range("I3").Value = 724.7895
$ws_BSM.Range("K3").Value = 724.7895
$ws_BSM.Range("M3").Value = -610.7895
$ws_BSM.Range("H20").Value = 2497.6365
$ws_BSM.Range("I20").Value = 2548.4
$ws_BSM.Range("K20").Value = 2548.4
$ws_BSM.Range("M20").Value = -2301.4
$ws_BSM.Range("H80").Value = 852.8333
$ws_BSM.Range("I80").Value = 323.75
$ws_BSM.Range("K80").Value = 323.75
$ws_BSM.Range("M80").Value = 674.25
$ws_BSM.Range("H82").Value = 26531.572
$ws_BSM.Range("I82").Value = 15628.5
$ws_BSM.Range("K82").Value = 15628.5
$ws_BSM.Range("M82").Value = -15245.5
$ws_BSM.Range("H83").Value = 852.8333
$ws_BSM.Range("I83").Value = 323.75
$ws_BSM.Range("K83").Value = 1618.75
$ws_BSM.Range("M83").Value = 3373.25
$ws_BSM.Range("H85").Value = 26531.572
$ws_BSM.Range("I85").Value = 15628.5
$ws_BSM.Range("K85").Value = 15628.5
$ws_BSM.Range("M85").Value = -14302.5
$ws_BSM.Range("H86").Value = 3204.5833
$ws_BSM.Range("I86").Value = 3353.6
$ws_BSM.Range("K86").Value = 3353.6
$ws_BSM.Range("M86").Value = -2230.6
$ws_BSM.Range("H89").Value = 3204.5833
$ws_BSM.Range("I89").Value = 3353.6
$ws_BSM.Range("K89").Value = 16768
$ws_BSM.Range("M89").Value = -11152
$ws_BSM.Range("H94").Value = 7353226.5
$ws_BSM.Range("I94").Value = 8064809.5
$ws_BSM.Range("J94").Value = 203
$ws_BSM.Range("K94").Value = 8064809.5
$ws_BSM.Range("L94").Value = 203
$ws_BSM.Range("M94").Value = -8064358.5
$ws_BSM.Range("N94").Value = -1105
$ws_BSM.Range("H105").Value = 47620252
$ws_BSM.Range("I105").Value = 52632590
$ws_BSM.Range("K105").Value = 52632590
$ws_BSM.Range("M105").Value = -52630843

# ---- CRP ----
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H31").Value = 1359.9615
$ws_CRP.Range("I31").Value = 1330.2954
$ws_CRP.Range("J31").Value = 1523.125
$ws_CRP.Range("K31").Value = 1330.2954
$ws_CRP.Range("L31").Value = 1523.125
$ws_CRP.Range("M31").Value = -1035.2954
$ws_CRP.Range("N31").Value = -2113.125
$ws_CRP.Range("H34").Value = 1359.9615
$ws_CRP.Range("I34").Value = 1330.2954
$ws_CRP.Range("J34").Value = 1523.125
$ws_CRP.Range("K34").Value = 1330.2954
$ws_CRP.Range("L34").Value = 1523.125
$ws_CRP.Range("M34").Value = -1128.2954
$ws_CRP.Range("N34").Value = -1927.125
$ws_CRP.Range("H81").Value = 16500
$ws_CRP.Range("J81").Value = 16500
$ws_CRP.Range("L81").Value = 16500
$ws_CRP.Range("N81").Value = -18496
$ws_CRP.Range("H84").Value = 16500
$ws_CRP.Range("J84").Value = 16500
$ws_CRP.Range("L84").Value = 49500
$ws_CRP.Range("N84").Value = -59484
$ws_CRP.Range("H132").Value = 2501.1333
$ws_CRP.Range("I132").Value = 1842.375
$ws_CRP.Range("K132").Value = 5527.125
$ws_CRP.Range("M132").Value = -2997.125

# ---- CUL ----
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H3").Value = 15050
$ws_CUL.Range("I3").Value = 8460
$ws_CUL.Range("J3").Value = 18644.545
$ws_CUL.Range("K3").Value = 25380
$ws_CUL.Range("L3").Value = 55933.63499999999
$ws_CUL.Range("M3").Value = -25268
$ws_CUL.Range("N3").Value = -56157.63499999999
$ws_CUL.Range("H55").Value = 2133.7856
$ws_CUL.Range("J55").Value = 2133.7856
$ws_CUL.Range("L55").Value = 6401.3568
$ws_CUL.Range("N55").Value = -6755.3568

# ---- LTW ----
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H2").Value = 385200
$ws_LTW.Range("J2").Value = 359090.9
$ws_LTW.Range("L2").Value = 359090.9
$ws_LTW.Range("N2").Value = -359314.9

# ---- WVR ----
$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H87").Value = 0
$ws_WVR.Range("J87").Value = 0
$ws_WVR.Range("L87").Value = 0
$ws_WVR.Range("N87").ClearContents()
$ws_WVR.Range("H90").Value = 0
$ws_WVR.Range("J90").Value = 0
$ws_WVR.Range("L90").Value = 0
$ws_WVR.Range("N90").ClearContents()
